$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row right after the last data row (row 16), inheriting its
# formatting (style "Ruim" used throughout the skills table).
$ws.Rows.Item(17).Insert()

# New skill entry: ID 15, OpressTheWeak, DamageSkill(has effect), Mana 20, CoolDown 1
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "OpressTheWeak"
$ws.Range("C17").Value = "DamageSkill(has effect)"
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = 1

# Column C ("Type") now holds the longer string "DamageSkill(has effect)" -
# resize it to fit the new content, same as the author's original bestFit column.
$ws.Columns.Item(3).AutoFit()

# Leave the selection on the newly entered row, as Excel would after the edit.
$ws.Rows.Item(17).Select()

$wb.Save()
